$d = $word.ActiveDocument

# Helper: replace the first occurrence of $find (anywhere in the body) with
# $replace, keeping the replacement in its own run. Toggling Bold off/on
# around the mutation stops the engine from silently merging the edited
# run into an adjacent run that happens to share identical formatting
# (which otherwise happens here because the placeholder runs sit right
# next to other bold-only runs, e.g. the lone space before "periode").
function Replace-Placeholder($find, $replace) {
    $rng = $d.Content
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
    $rng.Font.Bold = 0
    $rng.Text = $replace
    $rng.Font.Bold = 1
    return $rng
}

# 1. {{ opdrachtgever }} -> {{ vacature.opdrachtgever }}
Replace-Placeholder "{{ opdrachtgever }}" "{{ vacature.opdrachtgever }}" | Out-Null

# 2. {{ periode }} -> {{ vacature.periode }}
Replace-Placeholder "{{ periode }}" "{{ vacature.periode }}" | Out-Null

# 3. {{ dagen_per_week }} -> {{ vacature.dagen_per_week }}, then append a new
#    bold run containing a single space right after it (same paragraph).
$rng = Replace-Placeholder "{{ dagen_per_week }}" "{{ vacature.dagen_per_week }}"
$rng.Collapse(0)
$rng.InsertAfter(" ")
$rng.Bold = 1

# 4. {{ locatie }} -> {{ vacature.locatie }}
Replace-Placeholder "{{ locatie }}" "{{ vacature.locatie }}" | Out-Null
